# NIT-9017918616.xlsx — "Actualiza base de datos EC y agrega parte 1 de
# nuevos estado de cuenta"
#
# The workbook lists, per "Periodo Mora", one detail row for the worker.
# This edit:
#   1. Adds a new trailing detail row for period "2508" (same worker /
#      amounts as the other rows), pushing the signature block down by
#      one row.
#   2. Re-orders the existing "Periodo Mora" values so the list now reads
#      ascending (2501..2507) instead of descending (2507..2501).
#   3. Refreshes the two summary figures: "VALOR MORA" total and the
#      "Cant. Periodos" count (now 8 periods instead of 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# --- 1. Insert a blank row at 23 (shifts the signature rows, old 27/28,
#        down to 28/29).
$ws.Rows.Item(23).Insert()

# Clone row 22's formatting (the row that currently carries the heavier
# "last row" bottom border) down onto the new row 23, so the new
# trailing row gets that look.
$ws.Range("B22:J22").Copy()
$ws.Range("B23:J23").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = 0

# Give old row 22 the plain interior-row formatting (matching row 21)
# now that it is no longer the last row in the table.
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J22").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = 0

# --- 2. Fill the new row 23 with the new period's detail line.
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1140887488"
$ws.Range("D23").Value = "ENRIQUE ANTONIO PUELLO ROMERO"
$ws.Range("E23").Value = "2508"
$ws.Range("F23").Value = 56940
$ws.Range("G23").Value = 1423500

# --- 3. Re-order "Periodo Mora" (column E) for rows 16-22 so the table
#        reads ascending 2501..2507 (was descending 2507..2501).
$ws.Range("E16").Value = "2501"
$ws.Range("E17").Value = "2502"
$ws.Range("E18").Value = "2503"
$ws.Range("E19").Value = "2504"
$ws.Range("E20").Value = "2505"
$ws.Range("E21").Value = "2506"
$ws.Range("E22").Value = "2507"

# --- 4. Refresh the summary figures.
$ws.Range("E11").Value = 455520
$ws.Range("F13").Value = 8
